$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the participant row: new NIM and new name
$ws.Range("B3").Value = 2341760196
$ws.Range("C3").Value = "Kemal S"

# Move the active selection from D2 to B2
$ws.Range("B2").Select()
